$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Leading apostrophe forces Excel to store the value as literal
    # text (matches t="inlineStr"/shared-string cells in the source),
    # even when the text looks like a number (e.g. '301.24').
    $cell.Value = "'" + $text
    # Re-normalize so the quote-prefix / text-number-format tweak Excel
    # applies under the hood does not leave a stray cell style behind.
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "41.277.35"
Set-TextCell 2 5 "  +3.41%  "
Set-TextCell 3 4 "2.251.16"
Set-TextCell 3 5 "  +1.88%  "
Set-TextCell 4 5 "  +0.01%  "
Set-TextCell 5 4 "301.24"
Set-TextCell 5 5 "  +2.22%  "
Set-TextCell 6 4 "91.24"
Set-TextCell 6 5 "  +4.65%  "
Set-TextCell 7 4 "0.520"
Set-TextCell 7 5 "  +2.12%  "
Set-TextCell 8 5 "  -0.02%  "
Set-TextCell 9 4 "0.484"
Set-TextCell 9 5 "  +1.82%  "
Set-TextCell 10 2 "OKB"
Set-TextCell 10 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 10 4 "53.46"
Set-TextCell 10 5 "  +8.44%  "
Set-TextCell 11 2 "Avalanche"
Set-TextCell 11 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 11 4 "31.80"
Set-TextCell 11 5 "  +6.62%  "
Set-TextCell 12 4 "0.0793"
Set-TextCell 12 5 "  +2.19%  "
Set-TextCell 13 5 "  +2.95%  "
Set-TextCell 14 4 "6.58"
Set-TextCell 14 5 "  +1.98%  "
Set-TextCell 15 4 "2.604.05"
Set-TextCell 15 5 "  +1.88%  "
Set-TextCell 16 4 "14.10"
Set-TextCell 16 5 "  +2.58%  "
Set-TextCell 17 4 "2.298.62"
Set-TextCell 17 5 "  +3.60%  "
Set-TextCell 18 4 "0.750"
Set-TextCell 18 5 "  +2.79%  "
Set-TextCell 19 4 "41.225.78"
Set-TextCell 19 5 "  +3.49%  "
Set-TextCell 20 4 "11.86"
Set-TextCell 20 5 "  +5.65%  "
Set-TextCell 21 4 "0.0₃0901"
Set-TextCell 21 5 "  +1.95%  "
Set-TextCell 22 5 "  +1.34%  "
Set-TextCell 23 4 "66.72"
Set-TextCell 24 4 "239.66"
Set-TextCell 24 5 "  +1.74%  "
Set-TextCell 25 4 "2.55"
Set-TextCell 25 5 "  +3.56%  "
Set-TextCell 26 5 "  -0.01%  "
Set-TextCell 27 4 "1.86"
Set-TextCell 27 5 "  +3.43%  "
Set-TextCell 28 4 "23.70"
Set-TextCell 28 5 "  +5.19%  "
Set-TextCell 29 5 "  +4.32%  "
Set-TextCell 30 5 "  +4.07%  "
Set-TextCell 31 4 "159.52"
Set-TextCell 31 5 "  +2.44%  "
Set-TextCell 32 4 "33.41"
Set-TextCell 32 5 "  +4.99%  "
Set-TextCell 33 4 "1.00"
Set-TextCell 33 5 "  +0.19%  "
Set-TextCell 34 5 "  +5.69%  "
Set-TextCell 35 4 "0.0732"
Set-TextCell 35 5 "  +3.04%  "
Set-TextCell 36 4 "3.01"
Set-TextCell 36 5 "  +5.37%  "
Set-TextCell 37 4 "2.36"
Set-TextCell 37 5 "  +0.12%  "
Set-TextCell 38 4 "0.104"
Set-TextCell 38 5 "  +8.19%  "
Set-TextCell 39 4 "16.45"
Set-TextCell 39 5 "  +6.54%  "
Set-TextCell 40 5 "  +2.86%  "
Set-TextCell 41 4 "1.79"
Set-TextCell 41 5 "  +7.40%  "
Set-TextCell 42 5 "  +4.54%  "
Set-TextCell 43 4 "2.078.45"
Set-TextCell 43 5 "  -2.09%  "
Set-TextCell 44 4 "20.06"
Set-TextCell 44 5 "  +13.10%  "
Set-TextCell 45 4 "0.0277"
Set-TextCell 45 5 "  +4.00%  "
Set-TextCell 46 4 "10.21"
Set-TextCell 46 5 "  +5.08%  "
Set-TextCell 47 4 "2.97"
Set-TextCell 47 5 "  +11.58%  "
Set-TextCell 48 5 "  -3.10%  "
Set-TextCell 49 4 "2.473.27"
Set-TextCell 49 5 "  +2.00%  "
Set-TextCell 50 2 "TrustWalletToken"
Set-TextCell 50 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 50 4 "1.15"
Set-TextCell 50 5 "  +5.09%  "
Set-TextCell 51 2 "Stacks"
Set-TextCell 51 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell 51 4 "1.51"
Set-TextCell 51 5 "  +3.33%  "
